$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value (as scraped from coinranking.com by the GitHub Action).
# A leading apostrophe forces Excel to treat the assignment as literal text,
# which keeps number-looking strings (e.g. "7.50", "0.500") intact instead of
# being normalised into floats and losing trailing zeros / separators.
$updates = [ordered]@{
    "D2" = '26.789.97'
    "E2" = '  +0.00%  '
    "D3" = '1.638.14'
    "E3" = '  -0.43%  '
    "E4" = '  -0.47%  '
    "D5" = '218.78'
    "E5" = '  +0.72%  '
    "D6" = '0.500'
    "E6" = '  -0.53%  '
    "E7" = '  -0.32%  '
    "E8" = '  -0.74%  '
    "E9" = '  -0.84%  '
    "D10" = '19.18'
    "E10" = '  -0.05%  '
    "D11" = '0.0845'
    "E11" = '  +0.32%  '
    "D12" = '1.866.78'
    "E12" = '  -0.34%  '
    "D13" = '1.637.93'
    "E13" = '  -0.41%  '
    "E14" = '  -1.26%  '
    "D15" = '0.524'
    "E15" = '  -0.66%  '
    "D16" = '64.69'
    "E16" = '  +0.11%  '
    "D17" = '26.790.40'
    "E17" = '  -0.02%  '
    "D18" = '0.0₃0732'
    "E18" = '  -0.85%  '
    "D19" = '214.65'
    "E19" = '  +0.06%  '
    "E20" = '  -0.39%  '
    "E21" = '  -0.19%  '
    "D22" = '6.33'
    "E22" = '  +0.58%  '
    "D23" = '2.35'
    "E23" = '  -2.54%  '
    "D24" = '9.10'
    "E24" = '  -2.81%  '
    "D25" = '147.95'
    "E25" = '  +2.13%  '
    "E26" = '  -0.23%  '
    "E27" = '  -0.19%  '
    "E28" = '  -1.46%  '
    "D29" = '15.68'
    "E29" = '  -0.11%  '
    "D30" = '0.0505'
    "E30" = '  -1.84%  '
    "D31" = '1.20'
    "E31" = '  +1.21%  '
    "E32" = '  +1.39%  '
    "E33" = '  -0.32%  '
    "E34" = '  +0.06%  '
    "D35" = '1.260.91'
    "E35" = '  -2.11%  '
    "E36" = '  +0.23%  '
    "E37" = '  -0.54%  '
    "D38" = '0.528'
    "E38" = '  -2.38%  '
    "D39" = '0.814'
    "E39" = '  -1.52%  '
    "E40" = '  -0.29%  '
    "D41" = '0.804'
    "E41" = '  -1.02%  '
    "D42" = '5.33'
    "E42" = '  -0.36%  '
    "D43" = '1.777.99'
    "E43" = '  -0.91%  '
    "E44" = '  -4.50%  '
    "D45" = '92.16'
    "E45" = '  +0.75%  '
    "E46" = '  -0.20%  '
    "E47" = '  -1.71%  '
    "D48" = '0.0515'
    "E48" = '  -0.79%  '
    "D49" = '0.0961'
    "E49" = '  -1.62%  '
    "B50" = 'EnergySwap'
    "C50" = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    "D50" = '7.50'
    "E50" = '  -2.33%  '
    "B51" = 'USDD'
    "C51" = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
    "D51" = '1.01'
    "E51" = '  -0.24%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = "'" + $updates[$ref]
    $ws.Range($ref).Style = "Normal"
}
